$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.015247382432802
$ws.Range("D2").Value = 1.022334123525267
$ws.Range("E2").Value = 1.016933692870021
$ws.Range("F2").Value = 1.031856290586472
$ws.Range("I2").Value = 1.028306089089873
$ws.Range("J2").Value = 1.020473989073448
$ws.Range("K2").Value = 1.025168811009132
$ws.Range("L2").Value = 1.019784385885734
$ws.Range("M2").Value = 1.034663200481764
$ws.Range("N2").Value = 1.010947478118969

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.016092405973018
$ws.Range("D3").Value = 1.022981784943471
$ws.Range("E3").Value = 1.017646979328246
$ws.Range("F3").Value = 1.034017308502867
$ws.Range("I3").Value = 1.028527679182484
$ws.Range("J3").Value = 1.020954495095058
$ws.Range("K3").Value = 1.025623550988762
$ws.Range("L3").Value = 1.020303413555465
$ws.Range("M3").Value = 1.036629251772402
$ws.Range("N3").Value = 1.011106947497183

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.016639410432084
$ws.Range("D4").Value = 1.023400517598828
$ws.Range("E4").Value = 1.018109119858229
$ws.Range("F4").Value = 1.03540940210123
$ws.Range("I4").Value = 1.028668800320863
$ws.Range("J4").Value = 1.021265012290299
$ws.Range("K4").Value = 1.025916768310309
$ws.Range("L4").Value = 1.020639186959376
$ws.Range("M4").Value = 1.037894861278783
$ws.Range("N4").Value = 1.011209972413646

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.016869422775242
$ws.Range("D5").Value = 1.023576469327599
$ws.Range("E5").Value = 1.018303545700416
$ws.Range("F5").Value = 1.035993179428833
$ws.Range("I5").Value = 1.028727587544695
$ws.Range("J5").Value = 1.021395457159122
$ws.Range("K5").Value = 1.026039790440283
$ws.Range("L5").Value = 1.020780328368545
$ws.Range("M5").Value = 1.038425383988475
$ws.Range("N5").Value = 1.011253245091017

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.016908045876056
$ws.Range("D6").Value = 1.023606007493454
$ws.Range("E6").Value = 1.018336198954867
$ws.Range("F6").Value = 1.036091113601168
$ws.Range("I6").Value = 1.028737426550235
$ws.Range("J6").Value = 1.021417353756864
$ws.Range("K6").Value = 1.026060431949162
$ws.Range("L6").Value = 1.020804025579256
$ws.Range("M6").Value = 1.038514371678169
$ws.Range("N6").Value = 1.01126050847379

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.01664248366586
$ws.Range("D7").Value = 1.023402869002868
$ws.Range("E7").Value = 1.018111717229128
$ws.Range("F7").Value = 1.035417208255153
$ws.Range("I7").Value = 1.028669587959034
$ws.Range("J7").Value = 1.021266755680784
$ws.Range("K7").Value = 1.025918413106265
$ws.Range("L7").Value = 1.020641072968059
$ws.Range("M7").Value = 1.0379019561618
$ws.Range("N7").Value = 1.011210550778666

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.015532916746044
$ws.Range("D8").Value = 1.022553075555094
$ws.Range("E8").Value = 1.017174627361715
$ws.Range("F8").Value = 1.032587926287793
$ws.Range("I8").Value = 1.028381445895815
$ws.Range("J8").Value = 1.020636461499791
$ws.Range("K8").Value = 1.02532270583182
$ws.Range("L8").Value = 1.019959808406142
$ws.Range("M8").Value = 1.035329012768968
$ws.Range("N8").Value = 1.011001405110897

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.013579404269832
$ws.Range("D9").Value = 1.021052978419377
$ws.Range("E9").Value = 1.015527959550039
$ws.Range("F9").Value = 1.027553183874125
$ws.Range("I9").Value = 1.027856303978905
$ws.Range("J9").Value = 1.019522730553897
$ws.Range("K9").Value = 1.02426509887212
$ws.Range("L9").Value = 1.018758798614278
$ws.Range("M9").Value = 1.030743564721409
$ws.Range("N9").Value = 1.01063162409713

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.012278217817976
$ws.Range("D10").Value = 1.020051138520598
$ws.Range("E10").Value = 1.014433322204956
$ws.Range("F10").Value = 1.024161535245591
$ws.Range("I10").Value = 1.027494408443692
$ws.Range("J10").Value = 1.018778187175678
$ws.Range("K10").Value = 1.023554705346
$ws.Range("L10").Value = 1.017957788068788
$ws.Range("M10").Value = 1.027649971433243
$ws.Range("N10").Value = 1.010384274915562

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.011715067068287
$ws.Range("D11").Value = 1.019616909999473
$ws.Range("E11").Value = 1.013960085950809
$ws.Range("F11").Value = 1.022684106340707
$ws.Range("I11").Value = 1.027334880618098
$ws.Range("J11").Value = 1.018455304985018
$ws.Range("K11").Value = 1.023245830259261
$ws.Range("L11").Value = 1.017610865435497
$ws.Range("M11").Value = 1.026301286404574
$ws.Range("N11").Value = 1.010276974197593

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.011505928758305
$ws.Range("D12").Value = 1.019455554183342
$ws.Range("E12").Value = 1.013784418085905
$ws.Range("F12").Value = 1.022133957121293
$ws.Range("I12").Value = 1.027275198320568
$ws.Range("J12").Value = 1.018335298526041
$ws.Range("K12").Value = 1.023130908983478
$ws.Range("L12").Value = 1.017481991158533
$ws.Range("M12").Value = 1.025798913695859
$ws.Range("N12").Value = 1.010237088377892

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.011550787780366
$ws.Range("D13").Value = 1.019490168452623
$ws.Range("E13").Value = 1.013822094309086
$ws.Range("F13").Value = 1.022252028489681
$ws.Range("I13").Value = 1.027288019714433
$ws.Range("J13").Value = 1.01836104365171
$ws.Range("K13").Value = 1.023155568653339
$ws.Range("L13").Value = 1.017509635657466
$ws.Range("M13").Value = 1.025906738817935
$ws.Range("N13").Value = 1.010245645358754

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.011697778794135
$ws.Range("D14").Value = 1.019603573584018
$ws.Range("E14").Value = 1.013945562878505
$ws.Range("F14").Value = 1.022638658946511
$ws.Range("I14").Value = 1.027329955975319
$ws.Range("J14").Value = 1.018445386716075
$ws.Range("K14").Value = 1.023236334736653
$ws.Range("L14").Value = 1.017600212886354
$ws.Range("M14").Value = 1.02625978912718
$ws.Range("N14").Value = 1.010273677826105

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.011788350239771
$ws.Range("D15").Value = 1.01967343774202
$ws.Range("E15").Value = 1.014021650951078
$ws.Range("F15").Value = 1.022876692443285
$ws.Range("I15").Value = 1.027355737699969
$ws.Range("J15").Value = 1.018497343502989
$ws.Range("K15").Value = 1.023286072018958
$ws.Range("L15").Value = 1.017656018953834
$ws.Range("M15").Value = 1.02647712694808
$ws.Range("N15").Value = 1.010290945638632

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.012315598013821
$ws.Range("D16").Value = 1.020079947868775
$ws.Range("E16").Value = 1.014464745170352
$ws.Range("F16").Value = 1.024259397984195
$ws.Range("I16").Value = 1.027504936082336
$ws.Range("J16").Value = 1.018799605478492
$ws.Range("K16").Value = 1.023575177567944
$ws.Range("L16").Value = 1.017980810520012
$ws.Range("M16").Value = 1.027739283442664
$ws.Range("N16").Value = 1.010391391964921

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.012646399435196
$ws.Range("D17").Value = 1.02033482701089
$ws.Range("E17").Value = 1.014742887433655
$ws.Range("F17").Value = 1.025124341966592
$ws.Range("I17").Value = 1.027597766516898
$ws.Range("J17").Value = 1.018989075237982
$ws.Range("K17").Value = 1.023756185475695
$ws.Range("L17").Value = 1.018184522493062
$ws.Range("M17").Value = 1.028528527529028
$ws.Range("N17").Value = 1.010454346598969

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.01283937650267
$ws.Range("D18").Value = 1.020483452685212
$ws.Range("E18").Value = 1.014905195331953
$ws.Range("F18").Value = 1.025628000693771
$ws.Range("I18").Value = 1.027651640550235
$ws.Range("J18").Value = 1.019099542469492
$ws.Range("K18").Value = 1.023861641826771
$ws.Range("L18").Value = 1.018303336529967
$ws.Range("M18").Value = 1.028988001102647
$ws.Range("N18").Value = 1.010491047973967

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.012905181133329
$ws.Range("D19").Value = 1.020534123244727
$ws.Range("E19").Value = 1.014960550360991
$ws.Range("F19").Value = 1.025799592614542
$ws.Range("I19").Value = 1.02766996406659
$ws.Range("J19").Value = 1.01913720092217
$ws.Range("K19").Value = 1.023897578953834
$ws.Range("L19").Value = 1.018343847730852
$ws.Range("M19").Value = 1.029144521739542
$ws.Range("N19").Value = 1.010503558971004

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.012610904862904
$ws.Range("D20").Value = 1.020307485134607
$ws.Range("E20").Value = 1.014713037926391
$ws.Range("F20").Value = 1.025031629717029
$ws.Range("I20").Value = 1.027587834879952
$ws.Range("J20").Value = 1.018968751815846
$ws.Range("K20").Value = 1.023736777714612
$ws.Range("L20").Value = 1.018162666922247
$ws.Range("M20").Value = 1.028443940315757
$ws.Range("N20").Value = 1.010447594127414

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.011654492472644
$ws.Range("D21").Value = 1.019570180367543
$ws.Range("E21").Value = 1.013909201333735
$ws.Range("F21").Value = 1.022524843962841
$ws.Range("I21").Value = 1.027317618584711
$ws.Range("J21").Value = 1.018420551818858
$ws.Range("K21").Value = 1.023212556423727
$ws.Range("L21").Value = 1.017573540473113
$ws.Range("M21").Value = 1.026155863871265
$ws.Range("N21").Value = 1.01026542377888

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.011053395038713
$ws.Range("D22").Value = 1.019106237165387
$ws.Range("E22").Value = 1.013404452311371
$ws.Range("F22").Value = 1.020940802259705
$ws.Range("I22").Value = 1.027145253803936
$ws.Range("J22").Value = 1.018075450651263
$ws.Range("K22").Value = 1.022881850513488
$ws.Range("L22").Value = 1.017203065319515
$ws.Range("M22").Value = 1.024709076201609
$ws.Range("N22").Value = 1.010150715057229

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.011372025652624
$ws.Range("D23").Value = 1.019352217399612
$ws.Range("E23").Value = 1.013671967032815
$ws.Range("F23").Value = 1.021781297861089
$ws.Range("I23").Value = 1.027236862419771
$ws.Range("J23").Value = 1.018258435710095
$ws.Range("K23").Value = 1.023057269102706
$ws.Range("L23").Value = 1.017399467583127
$ws.Range("M23").Value = 1.025476834279015
$ws.Range("N23").Value = 1.010211540523017

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.012626943251823
$ws.Range("D24").Value = 1.020319839878701
$ws.Range("E24").Value = 1.014726525411138
$ws.Range("F24").Value = 1.025073525011963
$ws.Range("I24").Value = 1.02759232340091
$ws.Range("J24").Value = 1.018977935241863
$ws.Range("K24").Value = 1.02374554762523
$ws.Range("L24").Value = 1.01817254253808
$ws.Range("M24").Value = 1.028482164351458
$ws.Range("N24").Value = 1.010450645337456

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.01408423194487
$ws.Range("D25").Value = 1.021441102825134
$ws.Range("E25").Value = 1.015953112237939
$ws.Range("F25").Value = 1.028860826575824
$ws.Range("I25").Value = 1.027994138010908
$ws.Range("J25").Value = 1.019811019696332
$ws.Range("K25").Value = 1.024539452516558
$ws.Range("L25").Value = 1.019069349427853
$ws.Range("M25").Value = 1.03193532459825
$ws.Range("N25").Value = 1.010727367714823
